# Fix synapse IDs for CRC data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the synapse_id (column D) values for the CRC v2.0-public rows (71-82)
$ws.Range("D71").Value = "syn39802574"
$ws.Range("D72").Value = "syn39802575"
$ws.Range("D73").Value = "syn39802578"
$ws.Range("D74").Value = "syn39802579"
$ws.Range("D75").Value = "syn39802583"
$ws.Range("D76").Value = "syn39802586"
$ws.Range("D77").Value = "syn39802589"
$ws.Range("D78").Value = "syn39802591"
$ws.Range("D79").Value = "syn39802593"
$ws.Range("D80").Value = "syn39802604"
$ws.Range("D81").Value = "syn39802596"
$ws.Range("D82").Value = "syn39802627"

# Apply an AutoFilter across the full table range
$ws.Range("A1:E158").AutoFilter()

# Make sure the _FilterDatabase defined name is scoped to the sheet and hidden
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$158")
$n.Visible = $false

# Restore the view: scroll so row 105 is no longer pinned at top and select G81
$ws.Range("G81").Select()

$wb.Save()
